$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sub-header row (row 2: "(m3/s)" / "(MW)" / "(GWh)" labels).
# This shifts all data rows up by one (old row 3 -> new row 2, ... old row 11 -> new row 10).
$ws.Rows("2:2").Delete()

# Rewrite row 1 as a single header row with new column names.
# A1:E1 get the plain default style (the old E1 carried leftover formatting
# from the cell it replaces, so reset it back to Normal first).
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up a plain "font applied only" style (Arial 9, General format) distinct
# from the default style used by A1:E1. Create a transient named style to get the
# interop to mint that exact cellXf, then drop the named style again so the
# cellStyles/cellStyleXfs tables stay at their original size.
$tmpStyle = $wb.Styles.Add("TmpHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles.Item("TmpHeaderStyle").Delete()

# Move the selection to the first data row, matching the saved view state.
$ws.Range("A2:K2").Select()

Write-Output "done"
